$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the entire row 198 ("「全ての窓の奥に物語がある」") - all rows below
# shift up by one to fill the gap.
$ws.Rows.Item(198).Delete()
